# Izjednacavanje data excel fajlova
# Adds 5 new test rows (29-33) to Sheet1, extends the AutoFilter / filter
# database range accordingly, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-apply the AutoFilter over A1:F28 (the header row plus the 27 data
#     rows that existed before this edit) *before* the new rows are
#     appended below, so the engine doesn't auto-grow the filter range to
#     the eventual used range. Toggling off first avoids Range.AutoFilter()
#     flipping an already-existing filter off. ---
$ws.AutoFilterMode = $false
$ws.Range("A1:F28").AutoFilter()

# --- Keep the hidden _FilterDatabase defined name in sync with the new
#     AutoFilter range. ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$28"
    }
}

# --- Append the 5 new rows, cloning the formatting of the last existing
#     row (28) one at a time so every new row keeps its real values +
#     per-column styles (bulk multi-row paste only fills the first row). ---
$ws.Range("A28:F28").Copy($ws.Range("A29:F29"))
$ws.Range("A28:F28").Copy($ws.Range("A30:F30"))
$ws.Range("A28:F28").Copy($ws.Range("A31:F31"))
$ws.Range("A28:F28").Copy($ws.Range("A32:F32"))
$ws.Range("A28:F28").Copy($ws.Range("A33:F33"))

# Row 29
$ws.Cells.Item(29, 1).Value = "Manage_Products-Hide/Show_account_on_Product_List_[WEB]_1"
$ws.Cells.Item(29, 2).Value = "C70785"

# Row 30
$ws.Cells.Item(30, 1).Value = "Manage_Products-Hide/Show_account_on_Product_List_Invalid[WEB]_1"
$ws.Cells.Item(30, 2).Value = "C70786"

# Row 31
$ws.Cells.Item(31, 1).Value = "Manage_Products-Nickname_Product_Returning_To_Default_Name_[WEB]"
$ws.Cells.Item(31, 2).Value = "C70787"

# Row 32
$ws.Cells.Item(32, 1).Value = "Manage_Products-Nickname_Product_[WEB]_Invalid"
$ws.Cells.Item(32, 2).Value = "C70788"

# Row 33
$ws.Cells.Item(33, 1).Value = "Manage_Products-Nickname_Product_[WEB]"
$ws.Cells.Item(33, 2).Value = "C70789"

# --- Move the active selection / scroll position to match the edited
#     workbook (new last cell is B34, one below the appended data). ---
$ws.Cells.Item(15, 1).Select()
$ws.Range("B34").Select()
